# Update train schedule row 4 (x_nrSteps / y_nrSteps / alienID) and move the
# active selection to match, per the commit "updated task used in testing".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 data changes: D4 (x_corrSteps) 2 -> 3, F4 (y_nrSteps) 2 -> 3,
# H4 (alienID) 36 -> 46
$ws.Range("D4").Value = 3
$ws.Range("F4").Value = 3
$ws.Range("H4").Value = 46

# Move the active cell/selection from D5 to D4
$ws.Range("D4").Select()
